$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final roster data (player name, positions, team) that should occupy A2:C18
$data = @(
    @("Ja Morant", "PG", "Memphis Grizzlies"),
    @("Isaiah Collier", "PG,SG", "Utah Jazz"),
    @("De'Aaron Fox", "PG,SG", "San Antonio Spurs"),
    @("DeMar DeRozan", "SF,PF", "Sacramento Kings"),
    @("Christian Braun", "SG,SF", "Denver Nuggets"),
    @("Brook Lopez", "C", "Milwaukee Bucks"),
    @("Evan Mobley", "PF,C", "Cleveland Cavaliers"),
    @("Tyler Herro", "PG,SG", "Miami Heat"),
    @("Moussa Diabate", "C", "Charlotte Hornets"),
    @("Zach Collins", "PF,C", "Chicago Bulls"),
    @("Miles Bridges", "SF,PF", "Charlotte Hornets"),
    @("Luka Doncic", "PG,SG", "Los Angeles Lakers"),
    @("Josh Giddey", "PG,SG,SF", "Chicago Bulls"),
    @("Scottie Barnes", "PG,SG,SF,PF", "Toronto Raptors"),
    @("Mikal Bridges", "SG,SF,PF", "New York Knicks"),
    @("Nikola Vucevic", "PF,C", "Chicago Bulls"),
    @("P.J. Washington", "SF,PF", "Dallas Mavericks")
)

# Old sheet had 18 data rows (rows 2-19); new sheet only has 17 data rows (rows 2-18).
# Clear the previously used range first so the extra trailing row is removed.
$ws.Range("A2:C19").Clear()

$rowIndex = 2
foreach ($rec in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $rec[0]
    $ws.Cells.Item($rowIndex, 2).Value = $rec[1]
    $ws.Cells.Item($rowIndex, 3).Value = $rec[2]
    $rowIndex++
}
